$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7
$ws.Range("B7").Value = 17
$ws.Range("C7").Value = 24
$ws.Range("E7").Value = 8
$ws.Range("F7").Value = 15

# Row 8
$ws.Range("B8").Value = 23
$ws.Range("D8").Value = 7
$ws.Range("E8").Value = 14
$ws.Range("F8").Value = 16

# Row 9
$ws.Range("C9").Value = 6
$ws.Range("D9").Value = 13
$ws.Range("E9").Value = 20
$ws.Range("F9").Value = 22

# Row 10
$ws.Range("B10").Value = 10
$ws.Range("C10").Value = 12
$ws.Range("D10").Value = 19
$ws.Range("E10").Value = 21

# Row 11
$ws.Range("B11").Value = 11
$ws.Range("C11").Value = 18
$ws.Range("D11").Value = 25
$ws.Range("F11").Value = 9

# Update the selected cell to D11, matching the saved view state
$ws.Range("D11").Select()
